# ActionPlan.xlsx edit script
# Commit: "rewrite board, piece, pos"
#
# Adds three new "commit log" sections to the action-plan sheet:
#   - First implementation attempt inspired by Synesso/scala-chess (row 25, header only)
#   - Add credits (row 40 header + row 41 entry)
#   - Use ornicar.scalalib (row 44 header + rows 45-47 entries)
#   - Partial rewrite of Board, Piece and Pos (row 50 header + rows 51-53 entries)
# Also clears the stray fill styling left on A5 and A21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: drop the leftover "applyFill" style so it goes back to the default style ---
$ws.Range("A5").ClearFormats()

# --- Row 21: was highlighted in red; un-highlight it and add the matching "x" marker in B21 ---
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "x"

# --- New section header (row 25), inserted above the existing "Standard logic ..." (row 26) block ---
$ws.Range("A25").Value = "First implementation attempt inspired by Synesso/scala-chess"
$ws.Range("A25").Font.Bold = $true

# --- New section: "Add credits" ---
$ws.Range("A40").Value = "Add credits"
$ws.Range("A40").Font.Bold = $true

$ws.Range("A41").Value = "lila\README.md"
$ws.Range("B41").Value = "nr"

# --- New section: "Use ornicar.scalalib" ---
$ws.Range("A44").Value = "Use ornicar.scalalib"
$ws.Range("A44").Font.Bold = $true

$ws.Range("A45").Value = "lila\src\main\scala\package.scala"
$ws.Range("B45").Value = "nr"

$ws.Range("A46").Value = "lila\src\main\scala\Validation.scala"
$ws.Range("B46").Value = "x"

$ws.Range("A47").Value = "lila\project\Build.scala"
$ws.Range("B47").Value = "nr"

# --- New section: "Partial rewrite of Board, Piece and Pos" ---
$ws.Range("A50").Value = "Partial rewrite of Board, Piece and Pos"
$ws.Range("A50").Font.Bold = $true

$ws.Range("A51").Value = "lila\src\main\scala\model\Board.scala"
$ws.Range("B51").Value = "x"
$ws.Range("C51").Value = "For now - without Validation, but every function returns a new Board"

$ws.Range("A52").Value = "lila\src\main\scala\model\Piece.scala"
$ws.Range("B52").Value = "x"

$ws.Range("A53").Value = "lila\src\main\scala\model\Pos.scala"
$ws.Range("B53").Value = "x"

# --- Sheet view: keep selection/scroll roughly tracking the new bottom of the list ---
$ws.Range("A54").Select() | Out-Null
